# Progress log update: append a new "Week 8" section (rows 35-39) to the
# Hours worksheet, mirroring the existing weekly blocks above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week 8 header -------------------------------------------------------
$ws.Range("A35").Value = "Week 8"
$ws.Range("A35").Font.Bold = $true

# --- Day 1 (row 36) --------------------------------------------------------
$ws.Range("A36").Value = 42258
$ws.Range("A36").NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"

$ws.Range("B36").Value = 0.91666666666666663
$ws.Range("B36").NumberFormat = "h:mm AM/PM"

$ws.Range("C36").Value = 0
$ws.Range("C36").NumberFormat = "h:mm AM/PM"

$ws.Range("D36").Value = 0

$ws.Range("E36").Formula = "=MOD(C36-B36,1)*24-D36"
$ws.Range("E36").NumberFormat = "0.00"

# --- Day 2 (row 37) --------------------------------------------------------
$ws.Range("A37").Value = 42260
$ws.Range("A37").NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"

$ws.Range("B37").Value = 0.5
$ws.Range("B37").NumberFormat = "h:mm AM/PM"

$ws.Range("C37").Value = 0.041666666666666664
$ws.Range("C37").NumberFormat = "h:mm AM/PM"

$ws.Range("D37").Value = 2

$ws.Range("E37").Formula = "=MOD(C37-B37,1)*24-D37"
$ws.Range("E37").NumberFormat = "0.00"

# --- Day 3 (row 38) --------------------------------------------------------
$ws.Range("A38").Value = 42261
$ws.Range("A38").NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"

$ws.Range("B38").Value = 0.5
$ws.Range("B38").NumberFormat = "h:mm AM/PM"

$ws.Range("C38").Value = 0.041666666666666664
$ws.Range("C38").NumberFormat = "h:mm AM/PM"

$ws.Range("D38").Value = 1

$ws.Range("E38").Formula = "=MOD(C38-B38,1)*24-D38"
$ws.Range("E38").NumberFormat = "0.00"

# --- Week 8 total (row 39) --------------------------------------------------
$ws.Range("D39").Value = "Total"
$ws.Range("D39").Font.Bold = $true

$ws.Range("E39").Formula = "=SUM(E36:E38)"
$ws.Range("E39").NumberFormat = "0.00"
$ws.Range("E39").Font.Bold = $true

# --- Selection / view bookkeeping ------------------------------------------
[void]$ws.Range("C39").Select()
